# "Updated style and spacing"
#
# Trend_instructions (sheet 2):
#   - add a new column K header "override_normalization"
#   - widen column J (bestfit-ish) to fit its header text
#   - F5: 1 -> 2
#   - move the active-cell selection to F14
# Interpolation_instructions (sheet 1):
#   - move the active-cell selection to I1

$wb = $excel.ActiveWorkbook

$wsInterp = $wb.Worksheets.Item("Interpolation_instructions")
$wsTrend  = $wb.Worksheets.Item("Trend_instructions")

# --- Interpolation_instructions sheet -----------------------------------------

$wsInterp.Range("I1").Select()

# --- Trend_instructions sheet -------------------------------------------------

# New header in column K
$wsTrend.Range("K1").Value = "override_normalization"

# Give column J a bit more breathing room (manual "best fit" resize)
$wsTrend.Columns.Item(10).ColumnWidth = 13.25

# Update the data value in F5
$wsTrend.Range("F5").Value = 2

# Move the selection / active cell (also re-activates this sheet/tab last,
# matching the workbook's original active tab)
$wsTrend.Range("F14").Select()
